# Updated cryptos list on Fri Jul 19 15:43:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.539.23'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '3.453.52'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.32'
$ws.Range('E5').Value = '  +2.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.82'
$ws.Range('E6').Value = '  +6.68%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.452.53'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.562'
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.22'
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.430'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = '4.046.07'
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.51'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000176'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').Value = '65.474.82'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('D18').Value = '3.455.99'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.23'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.83'
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '382.82'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.93'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.66'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.521'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  +3.71%  '
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.32'
$ws.Range('E30').Value = '  +6.02%  '
$ws.Range('E31').Value = '  +3.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('E32').Value = '  +2.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.23'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.30'
$ws.Range('E34').Value = '  +5.64%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.52'
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.50'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +9.84%  '
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0741'
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.63'
$ws.Range('E41').Value = '  +5.50%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.807.53'
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.08'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.60'
$ws.Range('E44').Value = '  +5.06%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.03'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.46'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0310'
$ws.Range('E47').Value = '  +2.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.51'
$ws.Range('E48').Value = '  +6.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '348.68'
$ws.Range('E49').Value = '  +6.96%  '
$ws.Range('E50').Value = '  +2.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '32.28'
$ws.Range('E51').Value = '  +8.18%  '
